# Update gh-pages output data (苏州-漫展信息) to the values generated at 456a3b4.
# Two sheets ("展览" and "全部类型") carry the same per-event "want to go" (F)
# and "lowest ticket price" (G) columns that need refreshing.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (1st sheet) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 16382
$ws1.Range("F5").Value  = 422
$ws1.Range("F8").Value  = 15605
$ws1.Range("G9").Value  = 49.9
$ws1.Range("F10").Value = 9263
$ws1.Range("F11").Value = 476
$ws1.Range("F13").Value = 1030
$ws1.Range("F14").Value = 124
$ws1.Range("F15").Value = 221
$ws1.Range("F18").Value = 22
$ws1.Range("F20").Value = 614
$ws1.Range("F22").Value = 15
$ws1.Range("F23").Value = 78
$ws1.Range("F24").Value = 1154
$ws1.Range("F28").Value = 532
$ws1.Range("F34").Value = 63
$ws1.Range("F35").Value = 268
$ws1.Range("F36").Value = 366
$ws1.Range("F37").Value = 478
$ws1.Range("F39").Value = 5688
$ws1.Range("F40").Value = 5245

# ---- Sheet "全部类型" (4th sheet) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 16382
$ws4.Range("F5").Value  = 422
$ws4.Range("F8").Value  = 15605
$ws4.Range("G9").Value  = 49.9
$ws4.Range("F10").Value = 9263
$ws4.Range("F11").Value = 476
$ws4.Range("F13").Value = 1030
$ws4.Range("F14").Value = 124
$ws4.Range("F15").Value = 221
$ws4.Range("F18").Value = 22
$ws4.Range("F20").Value = 614
$ws4.Range("F22").Value = 15
$ws4.Range("F23").Value = 78
$ws4.Range("F24").Value = 1154
$ws4.Range("F28").Value = 532
$ws4.Range("F36").Value = 63
$ws4.Range("F37").Value = 268
$ws4.Range("F38").Value = 366
$ws4.Range("F39").Value = 478
$ws4.Range("F41").Value = 5688
$ws4.Range("F43").Value = 5245
